{"js": "const body = context.document.body;\n\n// ---------------------------------------------------------------------\n// 1) Team ID cell: \"PNT2022TMID4760\" -> \"PNT2022TMID47484\"\n//    (the author's edit landed as two runs: \"PNT2022TMID47\" + \"484\")\n// ---------------------------------------------------------------------\nconst teamIdResults = body.search(\"PNT2022TMID4760\", { matchCase: true });\nteamIdResults.load(\"items\");\nawait context.sync();\nteamIdResults.items[0].insertText(\"PNT2022TMID47484\", Word.InsertLocation.replace);\nawait context.sync();\n\n// Force the trailing \"484\" to land in its own run (matching run\n// formatting with its neighbour) by toggling bold on/off on just that\n// substring - this splits the run without changing any visible format.\nconst tailResults = body.search(\"484\", { matchCase: true });\ntailResults.load(\"items\");\nawait context.sync();\nconst tailRange = tailResults.items[0];\ntailRange.font.bold = true;\nawait context.sync();\ntailRange.font.bold = false;\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 2) Project Name cell: \"Child safety gadget using IoT\" ->\n//    \"IoT Based Safety Gadget for Child Safety Monitoring&Notification\"\n// ---------------------------------------------------------------------\nconst projectResults = body.search(\"Child safety gadget using IoT\", { matchCase: true });\nprojectResults.load(\"items\");\nawait context.sync();\nprojectResults.items[0].insertText(\n  \"IoT Based Safety Gadget for Child Safety Monitoring&Notification\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 3) Sprint-3/US-2 description cell: collapse the run-fragmented text\n//    (\"S\" + \"ervi\" + \"ce for accessing the IBM IoT platform.\") into a\n//    single run, and likewise for (\" \" + \"Create a Node-RED service.\").\n//    Re-inserting identical text over the whole matched span coalesces\n//    it into one run without touching the surrounding content.\n// ---------------------------------------------------------------------\nconst line1Results = body.search(\"Service for accessing the IBM IoT platform.\", { matchCase: true });\nline1Results.load(\"items\");\nawait context.sync();\nline1Results.items[0].insertText(\"Service for accessing the IBM IoT platform.\", Word.InsertLocation.replace);\nawait context.sync();\n\nconst line2Results = body.search(\" Create a Node-RED service.\", { matchCase: true });\nline2Results.load(\"items\");\nawait context.sync();\nline2Results.items[0].insertText(\" Create a Node-RED service.\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# 1) Team ID cell: \"PNT2022TMID4760\" -> \"PNT2022TMID47484\"\n#    (author's edit landed as two runs: \"PNT2022TMID47\" + \"484\")\n# ---------------------------------------------------------------------\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"PNT2022TMID4760\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"PNT2022TMID47484\"\n$find.Execute([ref]\"PNT2022TMID4760\", $false, $false, $false, $false, $false, $true, 1, $false, \"PNT2022TMID47484\", 2) | Out-Null\n\n# Force the replaced text to split into two runs (\"PNT2022TMID47\" / \"484\")\n# by toggling a character formatting property on the trailing \"484\" and\n# reverting it \u2014 this lands the new text as a separate run while leaving\n# the run formatting identical to its neighbour.\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"484\"\nif ($find2.Execute()) {\n    $sub = $d.Range($find2.Parent.Start, $find2.Parent.End)\n    $sub.Bold = 1\n    $sub.Bold = 0\n}\n\n# ---------------------------------------------------------------------\n# 2) Project Name cell: \"Child safety gadget using IoT\" ->\n#    \"IoT Based Safety Gadget for Child Safety Monitoring&Notification\"\n# ---------------------------------------------------------------------\n$find3 = $d.Content.Find\n$find3.ClearFormatting()\n$find3.Text = \"Child safety gadget using IoT\"\n$find3.Replacement.ClearFormatting()\n$find3.Replacement.Text = \"IoT Based Safety Gadget for Child Safety Monitoring&Notification\"\n$find3.Execute([ref]\"Child safety gadget using IoT\", $false, $false, $false, $false, $false, $true, 1, $false, \"IoT Based Safety Gadget for Child Safety Monitoring&Notification\", 2) | Out-Null\n\n# ---------------------------------------------------------------------\n# 3) Sprint-3/US-2 description cell: collapse the run-fragmented text\n#    (\"S\" + \"ervi\" + \"ce for accessing the IBM IoT platform.\") into a\n#    single run, and likewise for (\" \" + \"Create a Node-RED service.\").\n#    A Find/Replace-in-place (old text == new text) coalesces the\n#    matched span into one run without touching surrounding content.\n# ---------------------------------------------------------------------\n$find4 = $d.Content.Find\n$find4.ClearFormatting()\n$find4.Text = \"Service for accessing the IBM IoT platform.\"\n$find4.Replacement.ClearFormatting()\n$find4.Replacement.Text = \"Service for accessing the IBM IoT platform.\"\n$find4.Execute([ref]\"Service for accessing the IBM IoT platform.\", $false, $false, $false, $false, $false, $true, 1, $false, \"Service for accessing the IBM IoT platform.\", 2) | Out-Null\n\n$find5 = $d.Content.Find\n$find5.ClearFormatting()\n$find5.Text = \" Create a Node-RED service.\"\n$find5.Replacement.ClearFormatting()\n$find5.Replacement.Text = \" Create a Node-RED service.\"\n$find5.Execute([ref]\" Create a Node-RED service.\", $false, $false, $false, $false, $false, $true, 1, $false, \" Create a Node-RED service.\", 2) | Out-Null\n"}
